$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3307.6743
$ws.Range("J17").Value = 3307.6743
$ws.Range("L17").Value = 9923.0229
$ws.Range("N17").Value = -10259.0229

$ws.Range("H32").Value = 1040.6
$ws.Range("I32").Value = 601
$ws.Range("J32").Value = 1089.4445
$ws.Range("K32").Value = 601
$ws.Range("L32").Value = 1089.4445
$ws.Range("M32").Value = -275
$ws.Range("N32").Value = -1741.4445

$ws.Range("H33").Value = 34482984
$ws.Range("I33").Value = 37037264
$ws.Range("J33").Value = 175
$ws.Range("K33").Value = 37037264
$ws.Range("L33").Value = 175
$ws.Range("M33").Value = -37037035
$ws.Range("N33").Value = -633

$ws.Range("H125").Value = 2791.4443
$ws.Range("I125").Value = 4397.4
$ws.Range("J125").Value = 784
$ws.Range("K125").Value = 39576.6
$ws.Range("L125").Value = 7056
$ws.Range("M125").Value = -37116.6
$ws.Range("N125").Value = -11976

$ws.Range("H138").Value = 3273.158
$ws.Range("J138").Value = 3820.7
$ws.Range("L138").Value = 11462.1
$ws.Range("N138").Value = -21742.1

$ws.Range("H140").Value = 76483.336
$ws.Range("J140").Value = 76483.336
$ws.Range("L140").Value = 76483.336
$ws.Range("N140").Value = -86843.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 18520856
$ws.Range("I61").Value = 21741584
$ws.Range("J61").Value = 1678.5
$ws.Range("K61").Value = 21741584
$ws.Range("L61").Value = 1678.5
$ws.Range("M61").Value = -21741372
$ws.Range("N61").Value = -2102.5

$ws.Range("H132").Value = 7354997
$ws.Range("I132").Value = 11906606
$ws.Range("J132").Value = 2398.1538
$ws.Range("K132").Value = 35719818
$ws.Range("L132").Value = 7194.4614
$ws.Range("M132").Value = -35717288
$ws.Range("N132").Value = -12254.4614

$ws.Range("H136").Value = 18520856
$ws.Range("I136").Value = 21741584
$ws.Range("J136").Value = 1678.5
$ws.Range("K136").Value = 65224752
$ws.Range("L136").Value = 5035.5
$ws.Range("M136").Value = -65222202
$ws.Range("N136").Value = -10135.5

$ws.Range("H139").Value = 78130.71000000001
$ws.Range("J139").Value = 78130.71000000001
$ws.Range("L139").Value = 78130.71000000001
$ws.Range("N139").Value = -88410.71000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2900.6047
$ws.Range("I134").Value = 1646.4073
$ws.Range("K134").Value = 4939.2219
$ws.Range("M134").Value = -2404.2219

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8551751
$ws.Range("I31").Value = 7231.3335
$ws.Range("J31").Value = 18520358
$ws.Range("K31").Value = 7231.3335
$ws.Range("L31").Value = 18520358
$ws.Range("M31").Value = -6936.3335
$ws.Range("N31").Value = -18520948

$ws.Range("H34").Value = 8551751
$ws.Range("I34").Value = 7231.3335
$ws.Range("J34").Value = 18520358
$ws.Range("K34").Value = 7231.3335
$ws.Range("L34").Value = 18520358
$ws.Range("M34").Value = -7029.3335
$ws.Range("N34").Value = -18520762

$ws.Range("H58").Value = 1779.909
$ws.Range("I58").Value = 822.375
$ws.Range("J58").Value = 4333.3335
$ws.Range("K58").Value = 822.375
$ws.Range("L58").Value = 4333.3335
$ws.Range("M58").Value = -619.375
$ws.Range("N58").Value = -4739.3335

$ws.Range("H134").Value = 2042.6522
$ws.Range("I134").Value = 2059.842
$ws.Range("J134").Value = 1961
$ws.Range("K134").Value = 6179.526
$ws.Range("L134").Value = 5883
$ws.Range("M134").Value = -3644.526
$ws.Range("N134").Value = -10953

$ws.Range("H136").Value = 1779.909
$ws.Range("I136").Value = 822.375
$ws.Range("J136").Value = 4333.3335
$ws.Range("K136").Value = 2467.125
$ws.Range("L136").Value = 13000.0005
$ws.Range("M136").Value = 82.875
$ws.Range("N136").Value = -18100.0005

$ws.Range("H140").Value = 35600
$ws.Range("J140").Value = 35600
$ws.Range("L140").Value = 35600
$ws.Range("N140").Value = -45960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2013.35
$ws.Range("J22").Value = 2661.9285
$ws.Range("L22").Value = 7985.7855
$ws.Range("N22").Value = -8323.7855

$ws.Range("H27").Value = 2013.35
$ws.Range("J27").Value = 2661.9285
$ws.Range("L27").Value = 7985.7855
$ws.Range("N27").Value = -8189.7855

$ws.Range("H118").Value = 792.75
$ws.Range("J118").Value = 991.1177
$ws.Range("L118").Value = 2973.3531
$ws.Range("N118").Value = -5459.3531

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 44512.6
$ws.Range("I70").Value = 153245
$ws.Range("J70").Value = 4973.5454
$ws.Range("K70").Value = 153245
$ws.Range("L70").Value = 4973.5454
$ws.Range("M70").Value = -152975
$ws.Range("N70").Value = -5513.5454

$ws.Range("H73").Value = 44512.6
$ws.Range("I73").Value = 153245
$ws.Range("J73").Value = 4973.5454
$ws.Range("K73").Value = 153245
$ws.Range("L73").Value = 4973.5454
$ws.Range("M73").Value = -152309
$ws.Range("N73").Value = -6845.5454

$ws.Range("H132").Value = 4645.7427
$ws.Range("I132").Value = 3698.389
$ws.Range("J132").Value = 5648.8237
$ws.Range("K132").Value = 11095.167
$ws.Range("L132").Value = 16946.4711
$ws.Range("M132").Value = -8565.167000000001
$ws.Range("N132").Value = -22006.4711

$ws.Range("H138").Value = 57799.145
$ws.Range("J138").Value = 57799.145
$ws.Range("L138").Value = 57799.145
$ws.Range("N138").Value = -68079.14499999999

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1002.94116
$ws.Range("I46").Value = 907.1429000000001
$ws.Range("J46").Value = 1070
$ws.Range("K46").Value = 907.1429000000001
$ws.Range("L46").Value = 1070
$ws.Range("M46").Value = -719.1429000000001
$ws.Range("N46").Value = -1446

$ws.Range("H97").Value = 21644.5
$ws.Range("J97").Value = 21644.5
$ws.Range("L97").Value = 21644.5
$ws.Range("N97").Value = -23626.5

$ws.Range("H135").Value = 104982.25
$ws.Range("J135").Value = 104982.25
$ws.Range("L135").Value = 104982.25
$ws.Range("N135").Value = -115122.25

$ws.Range("H139").Value = 45054.89
$ws.Range("J139").Value = 45605.5
$ws.Range("L139").Value = 45605.5
$ws.Range("N139").Value = -55885.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2739.5
$ws.Range("I113").Value = 300
$ws.Range("J113").Value = 3349.375
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 10048.125
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -14388.125

$ws.Range("H132").Value = 1344.0566
$ws.Range("I132").Value = 1037.9025
$ws.Range("J132").Value = 2390.0833
$ws.Range("K132").Value = 3113.7075
$ws.Range("L132").Value = 7170.249899999999
$ws.Range("M132").Value = -583.7074999999995
$ws.Range("N132").Value = -12230.2499

$ws.Range("H136").Value = 1303.5385
$ws.Range("I136").Value = 1242.4762
$ws.Range("J136").Value = 1560
$ws.Range("K136").Value = 3727.4286
$ws.Range("L136").Value = 4680
$ws.Range("M136").Value = -1177.4286
$ws.Range("N136").Value = -9780

$ws.Range("H138").Value = 55167.5
$ws.Range("J138").Value = 65223.332
$ws.Range("L138").Value = 65223.332
$ws.Range("N138").Value = -75503.33199999999
